$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 6 ("te" / "adex489" / "2022-04-07" / <blank> / "fsd") down
# into the (previously unused) row 7, keeping the same literal text values
# and formatting as the source row.
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial()

# Row 6's "ID PRESTITO" (D6) blank placeholder is no longer needed.
$ws.Range("D6").ClearContents()

# Row 7 gets its own, new "ID PRESTITO" code.
$ws.Range("E7").Value = "nuovo codice"
